$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old stray formatted cell at A7 (no longer present after the edit) ---
$ws.Range("A7").Clear()

# --- New column widths for the newly-used D/E columns ---
$ws.Columns("D").ColumnWidth = 12.95
$ws.Columns("E").ColumnWidth = 9.6

# --- Row 6: Raunak / Naik / testing@example.com / 2123434565 / testing123 / test4567 ---
$ws.Range("A6").Value = "Raunak"
$ws.Range("B6").Value = "Naik"
$ws.Range("C6").Value = "testing@example.com"
$ws.Range("D6").Value = 2123434565
$ws.Range("E6").Value = "testing123"
$ws.Range("F6").Value = "test4567"

# --- Row 7: Raunak / Naik / demo1@example.com (hyperlinked) / 2123434565 / testing123 / testing123 ---
$ws.Range("A7").Value = "Raunak"
$ws.Range("B7").Value = "Naik"
$ws.Range("C7").Value = "demo1@example.com"
$ws.Range("D7").Value = 2123434565
$ws.Range("E7").Value = "testing123"
$ws.Range("F7").Value = "testing123"
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:demo1@example.com")

# --- Row 8: (blank first name) / Naik / testing1@example.com (hyperlinked) / 2123434565 / testing123 / testing123 ---
$ws.Range("B8").Value = "Naik"
$ws.Range("C8").Value = "testing1@example.com"
$ws.Range("D8").Value = 2123434565
$ws.Range("E8").Value = "testing123"
$ws.Range("F8").Value = "testing123"
$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:testing1@example.com")

# --- Row 9: Raunak / Naik / testing2@example.com (hyperlinked) / 12 / testing123 / testing123 ---
$ws.Range("A9").Value = "Raunak"
$ws.Range("B9").Value = "Naik"
$ws.Range("C9").Value = "testing2@example.com"
$ws.Range("D9").Value = 12
$ws.Range("E9").Value = "testing123"
$ws.Range("F9").Value = "testing123"
$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:testing2@example.com")

# --- Row 10: Raunakabcdefghijklmnopqrstuvwxyza / Raunak / testing3@example.com (hyperlinked) / 2123434565 / testing123 / testing123 ---
$ws.Range("C10").Value = "testing3@example.com"
$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:testing3@example.com")
$ws.Range("A10").Value = "Raunakabcdefghijklmnopqrstuvwxyza"
$ws.Range("B10").Value = "Raunak"
$ws.Range("D10").Value = 2123434565
$ws.Range("E10").Value = "testing123"
$ws.Range("F10").Value = "testing123"

# --- Row 11: Raunak / Naikabcdefghijklmnopqrstuvwxyzabc / testing4@example.com (hyperlinked) / 2123434565 / testing123 / testing123 ---
$ws.Range("A11").Value = "Raunak"
$ws.Range("B11").Value = "Naikabcdefghijklmnopqrstuvwxyzabc"
$ws.Range("C11").Value = "testing4@example.com"
$ws.Range("D11").Value = 2123434565
$ws.Range("E11").Value = "testing123"
$ws.Range("F11").Value = "testing123"
$ws.Hyperlinks.Add($ws.Range("C11"), "mailto:testing4@example.com")

# --- Row 12: Raunak / Naik / testing5@example.com (hyperlinked) / 2.12343456521312E+35 / testing123 / testing123 ---
$ws.Range("A12").Value = "Raunak"
$ws.Range("B12").Value = "Naik"
$ws.Range("C12").Value = "testing5@example.com"
$ws.Range("D12").Value = [double]"2.12343456521312E+35"
$ws.Range("E12").Value = "testing123"
$ws.Range("F12").Value = "testing123"
$ws.Hyperlinks.Add($ws.Range("C12"), "mailto:testing5@example.com")

# --- Apply the black explicit-colour font used by the newly entered rows (A6:F12) ---
$ws.Range("A6:F12").Font.Color = 0

# --- Selection / view state left by the editor ---
$ws.Range("D16").Select()
